$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks numeric need NumberFormat "@" forced
# BEFORE assigning the value, otherwise Excel auto-converts them to a
# real number (changing cell type from Text to Number).

$ws.Range("D2").Value = "30.542.80"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "2.137.39"
$ws.Range("E3").Value = "  +1.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.17"
$ws.Range("E5").Value = "  +5.32%  "
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5258"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4557"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.83"
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09163"
$ws.Range("E10").Value = "  +3.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.194"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.39"
$ws.Range("E12").Value = "  +5.39%  "
$ws.Range("D13").Value = "2.141.39"
$ws.Range("E13").Value = "  +2.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.899"
$ws.Range("E14").Value = "  +1.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.167"
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.02"
$ws.Range("E16").Value = "  +5.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001171"
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.007"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.43"
$ws.Range("E20").Value = "  +6.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.370"
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("D23").Value = "30.636.68"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.90"
$ws.Range("E24").Value = "  +4.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.385"
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("D26").Value = "2.390.52"
$ws.Range("E26").Value = "  +2.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.54"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.620"
$ws.Range("E28").Value = "  +4.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.81"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.04"
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.227"
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.717"
$ws.Range("E32").Value = "  +4.28%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.1082"
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.393"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.031"
$ws.Range("E35").Value = "  +2.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.135"
$ws.Range("E36").Value = "  +5.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.46"
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02650"
$ws.Range("E38").Value = "  +3.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06984"
$ws.Range("E39").Value = "  +2.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2343"
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.74"
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6987"
$ws.Range("E42").Value = "  +1.76%  "
$ws.Range("E43").Value = "  +2.54%  "
$ws.Range("E44").Value = "  +6.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.353"
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6516"
$ws.Range("E46").Value = "  +2.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000375"
$ws.Range("E47").Value = "  +8.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.751"
$ws.Range("E48").Value = "  +2.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.251"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.92"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07296"
$ws.Range("E51").Value = "  +2.37%  "
